# Add new columns I (I0) and J (IF) with per-row values, matching the
# source diff. Column I/J values for rows 2..36 (row 1 = headers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style used by the existing header cells (bold, bordered, centered)
# by copying the formatting from the last existing header cell (H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$values = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(7, 7)
    5  = @(8, 8)
    6  = @(10, 10)
    7  = @(7, 8)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(7, 8)
    11 = @(8, 8)
    12 = @(8, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(7, 7)
    16 = @(5, 6)
    17 = @(7, 8)
    18 = @(7, 8)
    19 = @(9, 9)
    20 = @(7, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(8, 9)
    25 = @(1, 4)
    26 = @(12, 12)
    27 = @(7, 7)
    28 = @(9, 9)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(6, 7)
    32 = @(6, 6)
    33 = @(9, 9)
    34 = @(6, 6)
    35 = @(9, 9)
    36 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
